# Update online UG model results stats (ug_envy) per commit: "online UG model results stats update"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 2.7668667708129
$ws.Range("C2").Value = 0.154719166873387
$ws.Range("D2").Value = 17.8831545355795
$ws.Range("E2").Value = 0.00000000000000000000000000000000000000000198838086898623

# Row 3 - depression_mc
$ws.Range("B3").Value = 0.130824900782017
$ws.Range("C3").Value = 0.313412707356774
$ws.Range("D3").Value = 0.417420537556865
$ws.Range("E3").Value = 0.676847136867132

# Row 4 - anhedonia_mc
$ws.Range("B4").Value = 0.0235433211371859
$ws.Range("C4").Value = 0.309449197675248
$ws.Range("D4").Value = 0.0760813772149233
$ws.Range("E4").Value = 0.939435257857624

# Row 5 - depression_mc:anhedonia_mc
$ws.Range("B5").Value = -0.0432914555763041
$ws.Range("C5").Value = 0.62561652674361
$ws.Range("D5").Value = -0.069198068985869
$ws.Range("E5").Value = 0.944905506194815
